$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.224.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.285.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.280.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.693.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.156.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.299.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("E20").Value = "  -2.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.51%  "

$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0724"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.79%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("E39").Value = "  -0.76%  "

$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "288.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("E50").Value = "  +0.26%  "

$ws.Range("E51").Value = "  +1.42%  "
